$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value2 = 0.2034883720930233
$ws.Cells.Item(2, 3).Value2 = 0.5494186046511628
$ws.Cells.Item(2, 10).Value2 = 0.008720930232558139
$ws.Cells.Item(2, 16).Value2 = 0.1424418604651163
$ws.Cells.Item(2, 19).Value2 = 0.09593023255813954
# Row 3
$ws.Cells.Item(3, 2).Value2 = 0.02116402116402116
$ws.Cells.Item(3, 3).Value2 = 0.01587301587301587
$ws.Cells.Item(3, 10).Value2 = 0.03174603174603174
$ws.Cells.Item(3, 16).Value2 = 0.708994708994709
$ws.Cells.Item(3, 19).Value2 = 0.2222222222222222
# Row 4
$ws.Cells.Item(4, 10).Value2 = 0.1
$ws.Cells.Item(4, 16).Value2 = 0.68
$ws.Cells.Item(4, 19).Value2 = 0.22
# Row 6
$ws.Cells.Item(6, 2).Value2 = 0.05284552845528456
$ws.Cells.Item(6, 4).Value2 = 0.01626016260162602
$ws.Cells.Item(6, 5).Value2 = 0.004065040650406504
$ws.Cells.Item(6, 6).Value2 = 0.05691056910569105
$ws.Cells.Item(6, 10).Value2 = 0.3739837398373984
$ws.Cells.Item(6, 15).Value2 = 0.01219512195121951
$ws.Cells.Item(6, 17).Value2 = 0.1016260162601626
$ws.Cells.Item(6, 18).Value2 = 0.08943089430894309
$ws.Cells.Item(6, 19).Value2 = 0.2926829268292683
# Row 7
$ws.Cells.Item(7, 2).Value2 = 0.1095238095238095
$ws.Cells.Item(7, 4).Value2 = 0.02380952380952381
$ws.Cells.Item(7, 6).Value2 = 0.05238095238095238
$ws.Cells.Item(7, 10).Value2 = 0.2095238095238095
$ws.Cells.Item(7, 15).Value2 = 0.009523809523809525
$ws.Cells.Item(7, 17).Value2 = 0.1476190476190476
$ws.Cells.Item(7, 18).Value2 = 0.07142857142857142
$ws.Cells.Item(7, 19).Value2 = 0.3761904761904762
# Row 8
$ws.Cells.Item(8, 2).Value2 = 0.0759493670886076
$ws.Cells.Item(8, 4).Value2 = 0.02025316455696203
$ws.Cells.Item(8, 5).Value2 = 0.002531645569620253
$ws.Cells.Item(8, 6).Value2 = 0.04556962025316456
$ws.Cells.Item(8, 10).Value2 = 0.1620253164556962
$ws.Cells.Item(8, 15).Value2 = 0.01265822784810127
$ws.Cells.Item(8, 17).Value2 = 0.1594936708860759
$ws.Cells.Item(8, 18).Value2 = 0.0810126582278481
$ws.Cells.Item(8, 19).Value2 = 0.4405063291139241
# Row 9
$ws.Cells.Item(9, 2).Value2 = 0.08374384236453201
$ws.Cells.Item(9, 4).Value2 = 0.009852216748768473
$ws.Cells.Item(9, 6).Value2 = 0.06403940886699508
$ws.Cells.Item(9, 10).Value2 = 0.167487684729064
$ws.Cells.Item(9, 15).Value2 = 0.01477832512315271
$ws.Cells.Item(9, 17).Value2 = 0.2167487684729064
$ws.Cells.Item(9, 18).Value2 = 0.0541871921182266
$ws.Cells.Item(9, 19).Value2 = 0.3891625615763547
# Row 10
$ws.Cells.Item(10, 2).Value2 = 0.1334776334776335
$ws.Cells.Item(10, 4).Value2 = 0.02453102453102453
$ws.Cells.Item(10, 5).Value2 = 0.001443001443001443
$ws.Cells.Item(10, 6).Value2 = 0.06926406926406926
$ws.Cells.Item(10, 10).Value2 = 0.1204906204906205
$ws.Cells.Item(10, 15).Value2 = 0.01298701298701299
$ws.Cells.Item(10, 17).Value2 = 0.1998556998556998
$ws.Cells.Item(10, 18).Value2 = 0.06277056277056277
$ws.Cells.Item(10, 19).Value2 = 0.3751803751803752
# Row 11
$ws.Cells.Item(11, 7).Value2 = 0.1405750798722045
$ws.Cells.Item(11, 10).Value2 = 0.1182108626198083
$ws.Cells.Item(11, 11).Value2 = 0.2108626198083067
$ws.Cells.Item(11, 12).Value2 = 0.5175718849840255
$ws.Cells.Item(11, 19).Value2 = 0.01277955271565495
# Row 12
$ws.Cells.Item(12, 7).Value2 = 0.7705882352941177
$ws.Cells.Item(12, 10).Value2 = 0.1647058823529412
$ws.Cells.Item(12, 11).Value2 = 0.005882352941176471
$ws.Cells.Item(12, 12).Value2 = 0.02941176470588235
$ws.Cells.Item(12, 19).Value2 = 0.02941176470588235
# Row 13
$ws.Cells.Item(13, 7).Value2 = 0.7580645161290323
$ws.Cells.Item(13, 10).Value2 = 0.1935483870967742
$ws.Cells.Item(13, 19).Value2 = 0.04838709677419355
# Row 15
$ws.Cells.Item(15, 6).Value2 = 0.01626016260162602
$ws.Cells.Item(15, 8).Value2 = 0.1869918699186992
$ws.Cells.Item(15, 9).Value2 = 0.08130081300813008
$ws.Cells.Item(15, 10).Value2 = 0.3739837398373984
$ws.Cells.Item(15, 11).Value2 = 0.06504065040650407
$ws.Cells.Item(15, 13).Value2 = 0.01219512195121951
$ws.Cells.Item(15, 15).Value2 = 0.04878048780487805
$ws.Cells.Item(15, 19).Value2 = 0.2154471544715447
# Row 16
$ws.Cells.Item(16, 6).Value2 = 0.03755868544600939
$ws.Cells.Item(16, 8).Value2 = 0.1267605633802817
$ws.Cells.Item(16, 9).Value2 = 0.0892018779342723
$ws.Cells.Item(16, 10).Value2 = 0.4037558685446009
$ws.Cells.Item(16, 11).Value2 = 0.136150234741784
$ws.Cells.Item(16, 13).Value2 = 0.02347417840375587
$ws.Cells.Item(16, 15).Value2 = 0.04225352112676056
$ws.Cells.Item(16, 19).Value2 = 0.1408450704225352
# Row 17
$ws.Cells.Item(17, 6).Value2 = 0.01834862385321101
$ws.Cells.Item(17, 8).Value2 = 0.1513761467889908
$ws.Cells.Item(17, 9).Value2 = 0.1123853211009174
$ws.Cells.Item(17, 10).Value2 = 0.3876146788990826
$ws.Cells.Item(17, 11).Value2 = 0.1077981651376147
$ws.Cells.Item(17, 13).Value2 = 0.02981651376146789
$ws.Cells.Item(17, 15).Value2 = 0.07339449541284404
$ws.Cells.Item(17, 19).Value2 = 0.1192660550458716
# Row 18
$ws.Cells.Item(18, 6).Value2 = 0.02395209580838323
$ws.Cells.Item(18, 8).Value2 = 0.1796407185628743
$ws.Cells.Item(18, 9).Value2 = 0.0658682634730539
$ws.Cells.Item(18, 10).Value2 = 0.4131736526946108
$ws.Cells.Item(18, 11).Value2 = 0.08383233532934131
$ws.Cells.Item(18, 13).Value2 = 0.01796407185628742
$ws.Cells.Item(18, 14).Value2 = 0.005988023952095809
$ws.Cells.Item(18, 15).Value2 = 0.09580838323353294
$ws.Cells.Item(18, 19).Value2 = 0.1137724550898204
# Row 19
$ws.Cells.Item(19, 6).Value2 = 0.03139356814701378
$ws.Cells.Item(19, 8).Value2 = 0.1753445635528331
$ws.Cells.Item(19, 9).Value2 = 0.07963246554364471
$ws.Cells.Item(19, 10).Value2 = 0.3790199081163859
$ws.Cells.Item(19, 11).Value2 = 0.1041347626339969
$ws.Cells.Item(19, 13).Value2 = 0.02909647779479326
$ws.Cells.Item(19, 15).Value2 = 0.09264931087289434
$ws.Cells.Item(19, 19).Value2 = 0.108728943338438
